# Updated Email Campaign - CF
$wb = $excel.ActiveWorkbook

# --- EmailLocationCampaign sheet: update location / feedback text ---
$ws1 = $wb.Worksheets.Item("EmailLocationCampaign")
$ws1.Activate()
$ws1.Range("F2").Value = "NTBACF02"
$ws1.Range("S2").Value = "Campaign Test06/12/2020 4:01:43 PM"
$ws1.Range("S2").Value = "Campaign Test06/12/2020 5:03:38 PM"
$ws1.Range("S2").Value = "Campaign Test06/12/2020 5:25:10 PM"
$ws1.Range("F2").Select()

# --- Date sheet: bump scheduled month references ---
$ws2 = $wb.Worksheets.Item("Date")
$ws2.Range("B3").Value = "August"
$ws2.Range("E3").Value = "September"

# --- Reschedule Date sheet: bump reschedule day/month ---
$ws3 = $wb.Worksheets.Item("Reschedule Date")
$ws3.Range("B3").Value = "August"
$ws3.Range("D3").Value = 30
# restore the quote-prefixed numeric style that plain Value writes clear
$ws3.Range("C3").Copy()
$ws3.Range("D3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws3.Range("E3").Value = "September"

# --- Add a new trailing worksheet carrying the location name ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newWs = $wb.Worksheets.Add($null, $lastSheet)
$newWs.Range("A1").Value = "Alexis Multispeciality Hospital"

# --- Reschedule Date stays the active/selected tab in the final state ---
$ws3.Activate()
$ws3.Range("E3").Select()
